# Remove the "password" column (D) from the user-import template.
# This mirrors removing the <c r="D1"> "password" shared string / cell
# and shifting the remaining columns (role, level) one place to the left,
# then leaves the active selection on J8 as in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D (the "password" header cell) and shift remaining cells left.
$ws.Range("D1").Delete(-4159) | Out-Null

# Match the final selection recorded in the saved workbook.
$ws.Range("J8").Select() | Out-Null
